# Commit: "Change the function name to start with __"
#
# 1. Rename worksheet "Fairs" -> "HistoricalFairs"
# 2. Rename the XLL function used in the array formula on the
#    (now) "HistoricalFairs" sheet from P_LoadReport to __LoadReport,
#    and tidy up the now-unused J:K column formatting / sheet dimension
#    that this leaves behind.
# 3. Update the "strategy" sheet's selection (view scrolled back to the
#    top, cursor moved to R7 instead of B23).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Fairs" sheet ---------------------------------------
$wsFairs = $wb.Worksheets.Item("Fairs")
$wsFairs.Name = "HistoricalFairs"

# --- 2. Update the report-loading formula on HistoricalFairs!A1 --------
$wsFairs.Range("A1:C784").FormulaArray = "=_xll.__LoadReport(1000, publish_config!B4,publish_config!H4,publish_config!I4)"

# Column I previously auto-sized to fit the old (wider) function output;
# shrink it back down now that the column only needs to hold numbers.
$wsFairs.Columns.Item(9).ColumnWidth = 9.333333333333334

# Columns J and K no longer carry any special formatting - clear it and
# drop the now-empty column definitions entirely.
$wsFairs.Columns.Item(10).ClearFormats()
$wsFairs.Columns.Item(11).ClearFormats()
$wsFairs.Columns.Item(10).Delete()
$wsFairs.Columns.Item(10).Delete()

# Move the selection on this sheet to I8.
$excel.Goto($wsFairs.Range("I8"))

# --- 3. Restore "strategy" as the active sheet / fix its selection -----
$wsStrategy = $wb.Worksheets.Item("strategy")
[void]$wsStrategy.Activate()
[void]$wsStrategy.Range("R7").Select()
